$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update rotor blade diameter row: label text drops "maximum " and value changes 265 -> 120
$ws.Range("A6").Value = "rotor blade diameter (m)"
$ws.Range("B6").Value = 120

# Keep active cell selection consistent with the target (B7)
$ws.Range("B7").Select()
